# Cambio de planes, Jesús
# - "Llenado de la BD" (row 8) is reassigned from Jesús to Sirio and its
#   dates are pushed back to the following week (matching the
#   "Vista de horarios de materias" slot dates).
# - "Documentación" (row 13) is reassigned from Jesús to Beatriz.
# - A new task "Manual de usuario" (row 14) is added, assigned to Beatriz.
# - The active cell selection moves from C8 to C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: "Llenado de la BD" -> responsible Sirio, dates moved a week later
$ws.Range("B8").Value = "Sirio"
$ws.Range("C8").Value = "Domingo 3 de mayo"
$ws.Range("D8").Value = "Domingo 3 de mayo"
$ws.Range("E8").Value = "Martes 5 de mayo"

# Row 13: "Documentación" -> responsible Beatriz
$ws.Range("B13").Value = "Beatriz"

# Row 14 (new): "Manual de usuario" -> responsible Beatriz
$ws.Range("A14").Value = "Manual de usuario"
$ws.Range("B14").Value = "Beatriz"

# Update the active selection to C13
$ws.Range("C13").Select()
